$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$red = 255
$yellow = 65535
$green = 5287936

# Fix typo in "Recomendations and next steps"
$ws.Range("A14").Value = "Recomendations and next steps"

# Status updates (A column fill color)
$ws.Range("A2").Interior.Color = $green    # Abstract: red -> green
$ws.Range("A4").Interior.Color = $green    # Customer Value Proposition: yellow -> green
$ws.Range("A5").Interior.Color = $green    # Changes from Proposal: yellow -> green
$ws.Range("A6").Interior.Color = $green    # Key Technical Elements: yellow -> green
$ws.Range("A8").Interior.Color = $green    # Failure Analysis: yellow -> green
$ws.Range("A9").Interior.Color = $green    # Timeline: yellow -> green
$ws.Range("A10").Interior.Color = $yellow  # Timeline changes writeup: red -> yellow
$ws.Range("A11").Interior.Color = $green   # Budget: yellow -> green
$ws.Range("A12").Interior.Color = $yellow  # Mid-review writeup: red -> yellow
$ws.Range("A13").Interior.Color = $yellow  # Process Understanding: red -> yellow
$ws.Range("A14").Interior.Color = $yellow  # Recomendations and next steps: red -> yellow
$ws.Range("A17").Interior.Color = $yellow  # Survey Responses: red -> yellow

$ws.Range("A11").Select()
